$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The previous export introduced an extra subtitle row ("population census
# results") under the title - remove it, restoring the original 2-row header.
$ws.Rows("2:2").Delete()

# The previous export also added stray 1989/2002 columns next to the 2014
# figures. Drop them so only the 2014 area value remains (fixing the export
# that was breaking the map generation).
$ws.Columns("B:C").Delete()

# Restore the proper sheet name (municipality name) instead of the generic "1".
$ws.Name = "ხაშური"

# Keep the saved cursor position/selection consistent with the cleaned sheet.
$ws.Range("A2").Select()
